$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B11").Value = "828959809"
$ws.Range("C11").Value = "3016877411"
$ws.Range("D11").Value = "732111198172294"

$ws.Range("B12").Value = "12669894"
$ws.Range("C12").Value = "3016876876"
$ws.Range("D12").Value = "732111198172293"

$ws.Rows.Item(13).Delete()

$ws.Range("C11").Select()
